$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.293.54'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -2.42%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.564.25'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -3.68%  '
$ws.Range('E4').Value = '  -0.32%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '206.93'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.15%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.01'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.29%  '
$ws.Range('E7').Value = '  -5.13%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.0605'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.63%  '
$ws.Range('E9').Value = '  -3.02%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '17.74'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.64%  '
$ws.Range('E11').Value = '  -0.81%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.784.30'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.50%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.571.78'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.25%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.00'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -4.34%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.504'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.67%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '25.290.03'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.41%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '59.24'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.14%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0₃0709'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -3.18%  '
$ws.Range('E19').Value = '  -0.39%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '184.96'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.32%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.12'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.48%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.26'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.14%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.87'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.19%  '
$ws.Range('E24').Value = '  -0.30%  '
$ws.Range('E25').Value = '  -4.40%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '139.84'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.64%  '
$ws.Range('E27').Value = '  -6.29%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.47'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.62%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '14.80'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.40%  '
$ws.Range('E30').Value = '  -6.43%  '
$ws.Range('E31').Value = '  -3.67%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.03'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.99%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.98'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.85%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.45'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.33%  '
$ws.Range('E35').Value = '  -3.59%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.088.00'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.69%  '
$ws.Range('E37').Value = '  -0.68%  '
$ws.Range('E38').Value = '  -4.94%  '
$ws.Range('E39').Value = '  -2.52%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.492'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -5.38%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.762'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -9.87%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.799'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +4.03%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '93.06'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -4.87%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.05'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.13%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.698.57'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.49%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0₆0108'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.65%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '52.42'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.70%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0504'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -4.61%  '
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.43'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.35%  '
$ws.Range('E50').Value = '  -1.69%  '
$ws.Range('E51').Value = '  -0.39%  '
